# Applies the 'Results with export bug fixed' update: appends 10 new
# result rows (15-24) to the Summary sheet, matching columns A:AF.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(46206429983.8214, 26009944462.88666, 0, 18165565760.93473, 2030919760, 44175510223.8214, 2030919760, 0, 0, 0, 0, 0, 46206429983.8214, 0, 0, 0, 0, 0, 0, 0, 0, 350.0540001392365, 46206422394.54113, 46206429983.8214, 7589.280265808105, 'optimal', 'costs_emissionlimit', -1, -1, 1, 'Island_lowDisc_greenfield', 'rawResults\20260215163629_Island_lowDisc_greenfield-1'),
    @(69075064785.38065, 49834925167.15388, 0, 17209219858.22677, 2030919760, 67044145025.38065, 2030919760, 0, 0, 0, 0, 0, 69075064785.38065, 0, 0, 0, 0, 0, 0, 0, 0, 583.933000087738, 69075063524.86021, 69075064785.38065, 1260.520431518555, 'optimal', 'costs_emissionlimit', -1, -1, 1, 'Island_highDisc_greenfield', 'rawResults\20260215174259_Island_highDisc_greenfield-1'),
    @(32379069972.6438, 21203768015.11133, 0, 15672035369.15194, 2030919760, 36875803384.26328, 2030919760, 1393576769.64288, -7921229941.262362, 0, 0, 0, 32379069972.6438, 0, 0, 0, 0, 0, 0, 0, 0, 337.835000038147, 32379068444.15659, 32379069972.6438, 1528.487209320068, 'optimal', 'costs_emissionlimit', -1, -1, 1, 'NotIsland_lowDisc_greenfield', 'rawResults\20260215185359_NotIsland_lowDisc_greenfield-1'),
    @(51104080452.09557, 39985471986.5044, 0, 14960541757.22156, 2030919760, 54946013743.72596, 2030919760, 1729431157.270354, -7602284208.900747, 0, 0, 0, 51104080452.09557, 0, 0, 0, 0, 0, 0, 0, 0, 999.9100000858307, 51104080209.30309, 51104080452.09557, 242.79248046875, 'optimal', 'costs_emissionlimit', -1, -1, 1, 'NotIsland_highDisc_greenfield', 'rawResults\20260215200052_NotIsland_highDisc_greenfield-1'),
    @(43455224245.01422, 24134714901.498, 0, 17289589583.51621, 2030919760, 41424304485.01422, 2030919760, 0, 0, 0, 0, 0, 43455224245.01422, 0, 0, 0, 0, 0, 0, 0, 0, 3202.686000108719, 43455223913.79698, 43455224245.01422, 331.2172393798828, 'optimal', 'costs_emissionlimit', -1, -1, 1, 'Island_lowDisc_greenfield', 'rawResults\20260216121620_Island_lowDisc_greenfield-1'),
    @(63105194110.90955, 44263260976.60344, 0, 16811013374.3061, 2030919760, 61074274350.90955, 2030919760, 0, 0, 0, 0, 0, 63105194110.90955, 0, 0, 0, 0, 0, 0, 0, 0, 4717.499000072479, 63105193966.09333, 63105194110.90955, 144.8162155151367, 'optimal', 'costs_emissionlimit', -1, -1, 1, 'Island_highDisc_greenfield', 'rawResults\20260216141800_Island_highDisc_greenfield-1'),
    @(43455224245.01422, 24134714901.498, 0, 17289589583.51621, 2030919760, 41424304485.01422, 2030919760, 0, 0, 0, 0, 0, 43455224245.01422, 0, 0, 0, 0, 0, 0, 0, 0, 3123.121000051498, 43455223913.79698, 43455224245.01422, 331.2172393798828, 'optimal', 'costs_emissionlimit', -1, -1, 1, 'Island_lowDisc_greenfield', 'rawResults\20260216200541_Island_lowDisc_greenfield-1'),
    @(63105194110.90955, 44263260976.60344, 0, 16811013374.3061, 2030919760, 61074274350.90955, 2030919760, 0, 0, 0, 0, 0, 63105194110.90955, 0, 0, 0, 0, 0, 0, 0, 0, 4508.098000049591, 63105193966.09333, 63105194110.90955, 144.8162155151367, 'optimal', 'costs_emissionlimit', -1, -1, 1, 'Island_highDisc_greenfield', 'rawResults\20260216220706_Island_highDisc_greenfield-1'),
    @(37433854532.10802, 19516430826.91516, 0, 14223510095.8617, 2030919760, 33739940922.77686, 2030919760, 1662993849.331164, 0, 0, 0, 0, 37433854532.10802, 0, 0, 0, 0, 0, 0, 0, 0, 1284.457999944687, 37433853089.32781, 37433854532.10802, 1442.780212402344, 'optimal', 'costs_emissionlimit', -1, -1, 1, 'NotIsland_lowDisc_greenfield', 'rawResults\20260217003106_NotIsland_lowDisc_greenfield-1'),
    @(53036906450.04081, 34764910461.50002, 0, 13789331371.06469, 2030919760, 48554241832.56471, 2030919760, 2451744857.476096, 0, 0, 0, 0, 53036906450.04081, 0, 0, 0, 0, 0, 0, 0, 0, 2715.637000083923, 53036905992.60873, 53036906450.04081, 457.4320755004883, 'optimal', 'costs_emissionlimit', -1, -1, 1, 'NotIsland_highDisc_greenfield', 'rawResults\20260217020104_NotIsland_highDisc_greenfield-1')
)

$startRow = 15
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowValues = $newRows[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowValues[$c]
    }
}
